$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.518.73"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -2.19%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.488.24"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.21%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.15"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.75%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "94.86"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -3.98%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -2.47%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -3.06%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "33.70"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -4.15%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -2.33%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.28%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.02"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -2.50%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.870.09"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.37%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.53"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.74%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.477.57"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.53%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.47%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "41.491.76"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -2.28%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -3.56%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0926"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.87%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.28"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -6.39%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "69.02"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.40%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.88"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.84%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -2.76%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.92"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -4.00%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.27"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -3.91%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.03%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -2.15%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.84"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -3.70%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "152.63"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.60%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.52"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -5.37%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -3.59%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "18.20"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +4.28%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -3.42%  "
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = "LidoDAOToken"
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.10"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.06%  "
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = "ApeXProtocol"
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.50"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -8.46%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.89"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -3.24%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.66%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -5.89%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.25"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +3.19%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.13%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "19.98"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -7.92%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.000.92"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.28%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -2.59%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.05"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -6.66%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.87"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.44%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.733.53"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.10%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "70.22"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.64%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "97.24"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -3.21%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -5.11%  "
